# Summary of the edit (per the diff):
#  - The "ssuid" values that lived in D2:F2 move down one row to D3:F3
#    (columns A:C on rows 2/3 are untouched).
#  - The sheet's frozen-pane top-left cell / active selection moves from
#    A17 / E29:E34 to A2 / D3:F3.
#  - The "duplicateValues" conditional-formatting rules that were scoped to
#    E2:E6 / F2:F6 are rescoped to E3:E6 / F3:F6 to track the moved data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the D2:F2 "ssuid" values down to D3:F3 --------------------------
# (NOTE: Range.Value is unreliable for string round-trips in this host --
# it surfaces a stringified member signature instead of the cell's actual
# content -- so Value2 is used here instead.)
$valD2 = $ws.Range("D2").Value2
$valE2 = $ws.Range("E2").Value2
$valF2 = $ws.Range("F2").Value2

$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()

$ws.Range("D3").Value2 = $valD2
$ws.Range("E3").Value2 = $valE2
$ws.Range("F3").Value2 = $valF2

# --- Rescope the conditional formatting that tracked the moved cells ------
$fcE = $ws.Range("E2:E6").FormatConditions.Item(1)
$fcE.ModifyAppliesToRange($ws.Range("E3:E6"))

$fcF = $ws.Range("F2:F6").FormatConditions.Item(1)
$fcF.ModifyAppliesToRange($ws.Range("F3:F6"))

# --- Update the view: frozen pane top-left cell + active selection --------
$aw = $ws.Application.ActiveWindow
$aw.ScrollRow = 2
$aw.ScrollColumn = 1
$ws.Range("D3:F3").Select()
